$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vpcheaders")
$ws.Activate()
$lo = $ws.ListObjects.Item(1)
Write-Host "before unlist: $($ws.ListObjects.Count)"
$lo.Unlist()
Write-Host "after unlist: $($ws.ListObjects.Count)"
$ws.Columns("B:B").Delete()
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:E3"), $null, 1)
Write-Host "after add: $($ws.ListObjects.Count)"
Write-Host "name: $($lo2.Name)"
